$wb = $excel.ActiveWorkbook

# --- Sheet "Appendix1" ---
# Header row (E1:G1) shifts to the new "paper" header strings (unchanged text,
# but re-pointed at new shared-string slots caused by the insert below) and a
# brand-new data row 2 is added for the paper's own author.
$ws1 = $wb.Worksheets.Item("Appendix1")
$ws1.Range("E1").Value = "Tên bài báo"
$ws1.Range("F1").Value = "Tên tạp chí"
$ws1.Range("G1").Value = "Ghi chú"

$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "Kieu Quoc Tuan"
$ws1.Range("C2").Value = "HE130005"
$ws1.Range("D2").Value = "FGWQN"
$ws1.Range("E2").Value = "paper 8"
$ws1.Range("F2").Value = "ABC"
$ws1.Range("G2").Value = "2 tác giả, 2 địa chỉ FPTU"

# --- Sheet "Apendix2" ---
# Existing paper (was "paper 4" / 1 author) becomes "paper 7" with 2 authors
# and 2 addresses; row 2 is updated in place and a new row 3 is appended for
# the second author.
$ws2 = $wb.Worksheets.Item("Apendix2")
$ws2.Range("E1").Value = "Tên bài báo"
$ws2.Range("F1").Value = "Tên tạp chí"
$ws2.Range("G1").Value = "Ghi chú"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "Kieu Quoc Tuan"
$ws2.Range("C2").Value = "HE130005"
$ws2.Range("D2").Value = "FGWQN"
$ws2.Range("E2").Value = "paper 7"
$ws2.Range("F2").Value = "ABC"
$ws2.Range("G2").Value = "2 tác giả, 2 địa chỉ FPTU"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Nguyễn Hồng Phúc"
$ws2.Range("C3").Value = "HE130001"
$ws2.Range("D3").Value = "FPTUHN2"
$ws2.Range("E3").Value = "paper 7"
$ws2.Range("F3").Value = "ABC"
$ws2.Range("G3").Value = "2 tác giả, 2 địa chỉ FPTU"

# --- Sheet "Apendix3" ---
# New payment-request row for the first author.
$ws3 = $wb.Worksheets.Item("Apendix3")
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "Kieu Quoc Tuan"
$ws3.Range("C2").Value = "HE130005"
$ws3.Range("D2").Value = "FGWQN"
$ws3.Range("E2").Value = "5.000.000 ₫"

# --- Sheet "Apendix4" ---
# Existing data row switches from "Nguyễn Văn A" to "Kieu Quoc Tuan" with the
# new amount, and a second row is appended for the second author.
$ws4 = $wb.Worksheets.Item("Apendix4")
$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = "Kieu Quoc Tuan"
$ws4.Range("C2").Value = "HE130005"
$ws4.Range("D2").Value = "FGWQN"
$ws4.Range("E2").Value = "5.000.000 ₫"

$ws4.Range("A3").Value = 2
$ws4.Range("B3").Value = "Nguyễn Hồng Phúc"
$ws4.Range("C3").Value = "HE130001"
$ws4.Range("D3").Value = "FPTUHN2"
$ws4.Range("E3").Value = "5.000.000 ₫"
